$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "risk/priority" legend value on the dashboard header (D1): Med-high -> Extreme
$ws.Range("D1").Value = "Extreme"

# Remove the two "Misc" rows (19 & 20) and blank out row 21 (Publications 2027 / IA program
# sections), leaving only row 21's first cell formatting behind - mirrors selecting A19:E21
# and deleting the contents while preserving A21's heading-row style.
$ws.Range("A22").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A19:E21").ClearContents()
$ws.Range("B19:E21").Clear()
$ws.Range("A19").Clear()
$ws.Range("A20:E20").Clear()
$ws.Rows(19).RowHeight = 15.05

# Move the active selection to D8 (was A5)
$ws.Range("D8").Select()
